# Applies the "Modified Absentees bug and consolidated output bug" change:
#  - Fill in the previously-empty rows 11-12 of the "Student Summary" sheet
#    with Course Code / EEOE 606 / Subject Code placeholder and
#    Max Marks / 40.
#  - Rename the "Average Marks & %" / bucket labels and recompute the
#    average-marks value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Student Summary")

# Rows 11-12 were previously an empty gap between the "Staff Name" row
# (row 10) and the "Attribute/Value" table (row 13 onward). Populate
# them with the new Course Code / Max Marks info instead of shifting
# any existing rows.
# New row 11: Course Code / EEOE 606 / <---- Type Subject Code
# New row 12: Max Marks / 40
# Copy the formatting (style) used by the row above (row 10) so the new
# rows match the rest of the info block.
$ws.Range("A10:D10").Copy()
$ws.Range("A11:D12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B11").Value = "Course Code:"
$ws.Range("C11").Value = "EEOE 606"
$ws.Range("D11").Value = "<---- Type Subject Code"

$ws.Range("B12").Value = "Max Marks:"
$ws.Range("C12").Value = 40

# Relabel the statistics bucket rows and refresh the computed average.
$ws.Range("A17").Value = "Average Marks"
$ws.Range("B17").Value = 22.22

$ws.Range("A18").Value = "Less Than 40%"
$ws.Range("A19").Value = "Between 40 % - 75 %"
$ws.Range("A20").Value = "More than 75%"
